# Add a new column W ("param_E_pv3_solar") to Sheet1, header in W1 and
# data values in W2:W11. Also update the existing U and V columns' data
# values to their new figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column W. Match the bold/bordered/centered header
# formatting used by the rest of row 1 (copy it from the neighbouring V1
# header cell), then set the text.
$ws.Range("V1").Copy()
$ws.Range("W1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("W1").Value = "param_E_pv3_solar"

# Updated U column values (rows 2-11).
$uValues = @(0.12, 0.12, 0.12, 0.12, 0.12, 0.12, 0.12, 0.12, 0.12, 0.12)

# Updated V column values (rows 2-11).
$vValues = @(49.93404166666668, 33.37400000000001, 157.15425, 0, 0, 0, 0, 0, 0, 0)

# New W column values (rows 2-11).
$wValues = @(33, 0, 0, 0, 0, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 21).Value = $uValues[$i]   # column U
    $ws.Cells.Item($row, 22).Value = $vValues[$i]   # column V
    $ws.Cells.Item($row, 23).Value = $wValues[$i]   # column W
}
